$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (shifts old rows 11-54 down to 12-55)
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new record
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = Get-Date -Year 2023 -Month 4 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100107
$ws.Cells.Item(11, 8).Value = "Otros"
$ws.Cells.Item(11, 9).Value = 100107001
$ws.Cells.Item(11, 10).Value = "Caqui"
$ws.Cells.Item(11, 11).Value = "Fuyu"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 80
$ws.Cells.Item(11, 14).Value = 24000
$ws.Cells.Item(11, 15).Value = 24000
$ws.Cells.Item(11, 16).Value = 24000
$ws.Cells.Item(11, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 1600
$ws.Cells.Item(11, 20).Value = 15
